$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '90.557.54'
Set-TextValue $ws.Range('E2') '  +1.57%  '

# Row 3
Set-TextValue $ws.Range('D3') '3.145.96'
Set-TextValue $ws.Range('E3') '  +4.11%  '

# Row 4
Set-TextValue $ws.Range('D4') '0.999'
Set-TextValue $ws.Range('E4') '  -0.18%  '

# Row 5
Set-TextValue $ws.Range('D5') '214.45'
Set-TextValue $ws.Range('E5') '  +2.49%  '

# Row 6
Set-TextValue $ws.Range('D6') '623.66'
Set-TextValue $ws.Range('E6') '  +2.07%  '

# Row 7
Set-TextValue $ws.Range('E7') '  +30.57%  '

# Row 8
Set-TextValue $ws.Range('D8') '0.367'
Set-TextValue $ws.Range('E8') '  +2.57%  '

# Row 9
Set-TextValue $ws.Range('D9') '1.00'
Set-TextValue $ws.Range('E9') '  -0.02%  '

# Row 10
Set-TextValue $ws.Range('D10') '3.144.36'
Set-TextValue $ws.Range('E10') '  +4.07%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.764'
Set-TextValue $ws.Range('E11') '  +15.39%  '

# Row 12
Set-TextValue $ws.Range('E12') '  +8.70%  '

# Row 13
Set-TextValue $ws.Range('B13') 'Toncoin'
Set-TextValue $ws.Range('C13') 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D13') '5.66'
Set-TextValue $ws.Range('E13') '  +6.46%  '

# Row 14
Set-TextValue $ws.Range('B14') 'ShibaInu'
Set-TextValue $ws.Range('C14') 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D14') '0.0000243'
Set-TextValue $ws.Range('E14') '  +3.51%  '

# Row 15
Set-TextValue $ws.Range('E15') '  +9.89%  '

# Row 16
Set-TextValue $ws.Range('D16') '90.305.27'
Set-TextValue $ws.Range('E16') '  +1.76%  '

# Row 17
Set-TextValue $ws.Range('D17') '3.732.14'
Set-TextValue $ws.Range('E17') '  +3.87%  '

# Row 18
Set-TextValue $ws.Range('D18') '3.190.34'
Set-TextValue $ws.Range('E18') '  +5.43%  '

# Row 19
Set-TextValue $ws.Range('D19') '3.64'
Set-TextValue $ws.Range('E19') '  +9.51%  '

# Row 20
Set-TextValue $ws.Range('D20') '14.24'
Set-TextValue $ws.Range('E20') '  +7.52%  '

# Row 21
Set-TextValue $ws.Range('D21') '462.52'
Set-TextValue $ws.Range('E21') '  +9.88%  '

# Row 22
Set-TextValue $ws.Range('D22') '0.0000209'
Set-TextValue $ws.Range('E22') '  -0.04%  '

# Row 23
Set-TextValue $ws.Range('D23') '9.06'
Set-TextValue $ws.Range('E23') '  +12.27%  '

# Row 24
Set-TextValue $ws.Range('D24') '5.28'
Set-TextValue $ws.Range('E24') '  +6.62%  '

# Row 25
Set-TextValue $ws.Range('D25') '5.82'
Set-TextValue $ws.Range('E25') '  +9.40%  '

# Row 26
Set-TextValue $ws.Range('D26') '89.29'
Set-TextValue $ws.Range('E26') '  +7.92%  '

# Row 27
Set-TextValue $ws.Range('D27') '12.00'
Set-TextValue $ws.Range('E27') '  +3.73%  '

# Row 28
Set-TextValue $ws.Range('D28') '3.313.20'
Set-TextValue $ws.Range('E28') '  +3.68%  '

# Row 29
Set-TextValue $ws.Range('E29') '  +0.01%  '

# Row 30
Set-TextValue $ws.Range('E30') '  +1.72%  '

# Row 31
Set-TextValue $ws.Range('D31') '0.161'
Set-TextValue $ws.Range('E31') '  +0.18%  '

# Row 32
Set-TextValue $ws.Range('D32') '9.17'
Set-TextValue $ws.Range('E32') '  +13.14%  '

# Row 33
Set-TextValue $ws.Range('D33') '27.11'
Set-TextValue $ws.Range('E33') '  +20.01%  '

# Row 34
Set-TextValue $ws.Range('D34') '515.52'
Set-TextValue $ws.Range('E34') '  +3.34%  '

# Row 35
Set-TextValue $ws.Range('D35') '0.184'
Set-TextValue $ws.Range('E35') '  +35.60%  '

# Row 36
Set-TextValue $ws.Range('D36') '3.60'
Set-TextValue $ws.Range('E36') '  +1.00%  '

# Row 37
Set-TextValue $ws.Range('E37') '  +6.90%  '

# Row 38
Set-TextValue $ws.Range('E38') '  +9.69%  '

# Row 39
Set-TextValue $ws.Range('D39') '6.84'
Set-TextValue $ws.Range('E39') '  +4.38%  '

# Row 40
Set-TextValue $ws.Range('E40') '  +5.04%  '

# Row 41
Set-TextValue $ws.Range('D41') '0.0871'
Set-TextValue $ws.Range('E41') '  +30.20%  '

# Row 42
Set-TextValue $ws.Range('D42') '22.19'
Set-TextValue $ws.Range('E42') '  -0.08%  '

# Row 43
Set-TextValue $ws.Range('E43') '  +0.26%  '

# Row 44
Set-TextValue $ws.Range('D44') '0.409'
Set-TextValue $ws.Range('E44') '  +14.26%  '

# Row 45
Set-TextValue $ws.Range('D45') '1.94'
Set-TextValue $ws.Range('E45') '  +7.41%  '

# Row 46
Set-TextValue $ws.Range('E46') '  +0.05%  '

# Row 47
Set-TextValue $ws.Range('D47') '4.59'
Set-TextValue $ws.Range('E47') '  +14.25%  '

# Row 48
Set-TextValue $ws.Range('B48') 'Monero'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D48') '148.62'
Set-TextValue $ws.Range('E48') '  +2.06%  '

# Row 49
Set-TextValue $ws.Range('B49') 'OKB'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D49') '45.33'
Set-TextValue $ws.Range('E49') '  +4.82%  '

# Row 50
Set-TextValue $ws.Range('E50') '  +11.41%  '

# Row 51
Set-TextValue $ws.Range('D51') '0.664'
Set-TextValue $ws.Range('E51') '  +14.81%  '
